$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: dylan
$ws.Range("A4").Value = "dylan"
$ws.Range("B4").Value = "dylan@email.com"
$ws.Range("C4").Value = "melbourne"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:dylan@email.com")
$ws.Range("B4").Style = "Hyperlink"

# Row 5: eric
$ws.Range("A5").Value = "eric"
$ws.Range("B5").Value = "eric@email.com"
$ws.Range("C5").Value = "brisbane"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:eric@email.com")
$ws.Range("B5").Style = "Hyperlink"

# Update the active selection to match the new last cell
$null = $ws.Range("C5").Select()
